# Progression script insertion donnees
# Adds idFormateur numbering to "Formateurs", and three new reference/association
# sheets: "Matieres", "Constitutions" and "Animations".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Formateurs: add an "idFormateur" id column (A) in front of the existing data
# ---------------------------------------------------------------------------
$wsFormateurs = $wb.Worksheets.Item("Formateurs")
$wsFormateurs.Range("A1").Value = "idFormateur"
$wsFormateurs.Range("A2").Value = 1
$wsFormateurs.Range("A3").Value = 2
$wsFormateurs.Range("A4").Value = 3
$wsFormateurs.Range("A5").Value = 4
$wsFormateurs.Columns.Item(1).AutoFit()

# ---------------------------------------------------------------------------
# 2. New sheet "Matieres" (placed right after "Formations")
# ---------------------------------------------------------------------------
$wsFormations = $wb.Worksheets.Item("Formations")
$wsMatieres = $wb.Worksheets.Add($null, $wsFormations)
$wsMatieres.Name = "Matieres"

$wsMatieres.Range("A1").Value = "idMatiere"
$wsMatieres.Range("B1").Value = "Matières"

$wsMatieres.Range("A2").Value = 1
$wsMatieres.Range("B2").Value = "Sport"
$wsMatieres.Range("C2").Formula = '="INSERT INTO Matieres (nomMatiere) VALUES (""" & B2 & """);"'

$wsMatieres.Range("A3").Value = 2
$wsMatieres.Range("B3").Value = "Français"
$wsMatieres.Range("C3").Formula = '="INSERT INTO Matieres (nomMatiere) VALUES (""" & B3 & """);"'

$wsMatieres.Range("A4").Value = 3
$wsMatieres.Range("B4").Value = "Math"
$wsMatieres.Range("C4").Formula = '="INSERT INTO Matieres (nomMatiere) VALUES (""" & B4 & """);"'

$wsMatieres.Columns.Item(3).AutoFit()

# ---------------------------------------------------------------------------
# 3. New sheet "Constitutions" (placed right after "Matieres")
#    cross-reference of Matieres x Formations
# ---------------------------------------------------------------------------
$wsConstitutions = $wb.Worksheets.Add($null, $wsMatieres)
$wsConstitutions.Name = "Constitutions"

$wsConstitutions.Range("A1").Value = "Matières"
$wsConstitutions.Range("B1").Value = "idMatiere"
$wsConstitutions.Range("B1:C1").Borders.Item(10).LineStyle = 1
$wsConstitutions.Range("B1:C1").Borders.Item(10).Weight = -4138
$wsConstitutions.Range("D1").Value = "idFormation"
$wsConstitutions.Range("E1").Value = "Formation"
$wsConstitutions.Range("F1").Value = "idMatiere"

$wsConstitutions.Range("A2").Value = "Sport"
$wsConstitutions.Range("B2").Value = 1
$wsConstitutions.Range("D2").Value = 1
$wsConstitutions.Range("E2").Value = "TSAII"
$wsConstitutions.Range("F2").Formula = '=VLOOKUP(B10,$A$2:$B$4,2,FALSE)'
$wsConstitutions.Range("G2").Formula = '="INSERT INTO constitutions (idMatiere, idFormation) VALUES ("&D2&","&F2&");"'

$wsConstitutions.Range("A3").Value = "Français"
$wsConstitutions.Range("B3").Value = 2
$wsConstitutions.Range("D3").Value = 2
$wsConstitutions.Range("E3").Value = "TRTE"
$wsConstitutions.Range("F3").Formula = '=VLOOKUP(B11,$A$2:$B$4,2,FALSE)'
$wsConstitutions.Range("G3").Formula = '="INSERT INTO constitutions (idMatiere, idFormation) VALUES ("&D3&","&F3&");"'

$wsConstitutions.Range("A4").Value = "Math"
$wsConstitutions.Range("B4").Value = 3
$wsConstitutions.Range("D4").Value = 3
$wsConstitutions.Range("E4").Value = "DWWM"
$wsConstitutions.Range("F4").Formula = '=VLOOKUP(B12,$A$2:$B$4,2,FALSE)'
$wsConstitutions.Range("G4").Formula = '="INSERT INTO constitutions (idMatiere, idFormation) VALUES ("&D4&","&F4&");"'

$wsConstitutions.Range("B2:C4").Borders.Item(10).LineStyle = 1
$wsConstitutions.Range("B2:C4").Borders.Item(10).Weight = 2
$wsConstitutions.Range("B2:C4").Locked = $true

$wsConstitutions.Range("D5").Value = 4
$wsConstitutions.Range("E5").Value = "CDA"
$wsConstitutions.Range("F5").Formula = '=VLOOKUP(B13,$A$2:$B$4,2,FALSE)'
$wsConstitutions.Range("G5").Formula = '="INSERT INTO constitutions (idMatiere, idFormation) VALUES ("&D5&","&F5&");"'

$wsConstitutions.Range("D6").Value = 5
$wsConstitutions.Range("E6").Value = "TSSR"
$wsConstitutions.Range("F6").Formula = '=VLOOKUP(B14,$A$2:$B$4,2,FALSE)'
$wsConstitutions.Range("G6").Formula = '="INSERT INTO constitutions (idMatiere, idFormation) VALUES ("&D6&","&F6&");"'

$wsConstitutions.Range("A9").Value = "Formation"
$wsConstitutions.Range("B9").Value = "matieres"

$wsConstitutions.Range("A10").Value = "TSAII"
$wsConstitutions.Range("B10").Value = "Sport"

$wsConstitutions.Range("A11").Value = "TRTE"
$wsConstitutions.Range("B11").Value = "Math"

$wsConstitutions.Range("A12").Value = "DWWM"
$wsConstitutions.Range("B12").Value = "Français"

$wsConstitutions.Range("A13").Value = "CDA"
$wsConstitutions.Range("B13").Value = "Français"

$wsConstitutions.Range("A14").Value = "TSSR"
$wsConstitutions.Range("B14").Value = "Math"

$wsConstitutions.Columns.Item(1).AutoFit()
$wsConstitutions.Columns.Item(4).AutoFit()
$wsConstitutions.Columns.Item(5).AutoFit()
$wsConstitutions.Columns.Item(7).AutoFit()

# ---------------------------------------------------------------------------
# 4. New sheet "Animations" (placed right after "Constitutions")
#    cross-reference of Formateurs x Formations
# ---------------------------------------------------------------------------
$wsAnimations = $wb.Worksheets.Add($null, $wsConstitutions)
$wsAnimations.Name = "Animations"

$wsAnimations.Range("A1").Value = "idFormateur"
$wsAnimations.Range("A1").Borders.Item(10).LineStyle = 1
$wsAnimations.Range("A1").Borders.Item(10).Weight = -4138
$wsAnimations.Range("B1").Value = "Formateur"

$wsAnimations.Range("A2").Value = 1
$wsAnimations.Range("B2").Value = "Poix"
$wsAnimations.Range("C2").Value = "Martine"

$wsAnimations.Range("A3").Value = 2
$wsAnimations.Range("B3").Value = "Dubois"
$wsAnimations.Range("C3").Value = "Thomas"

$wsAnimations.Range("A4").Value = 3
$wsAnimations.Range("B4").Value = "Butterdroghe"
$wsAnimations.Range("C4").Value = "Hervé"

$wsAnimations.Range("A5").Value = 4
$wsAnimations.Range("B5").Value = "Batzic"
$wsAnimations.Range("C5").Value = "Jean Paul"

$wsAnimations.Range("A2:A5").Borders.Item(10).LineStyle = 1
$wsAnimations.Range("A2:A5").Borders.Item(10).Weight = 2
$wsAnimations.Range("A2:A5").Locked = $true

$wsAnimations.Range("A7").Value = "idFormation"
$wsAnimations.Range("B7").Value = "Formation"

$wsAnimations.Range("A8").Value = 1
$wsAnimations.Range("B8").Value = "TSAII"

$wsAnimations.Range("A9").Value = 2
$wsAnimations.Range("B9").Value = "TRTE"

$wsAnimations.Range("A10").Value = 3
$wsAnimations.Range("B10").Value = "DWWM"

$wsAnimations.Range("A11").Value = 4
$wsAnimations.Range("B11").Value = "CDA"

$wsAnimations.Range("A12").Value = 5
$wsAnimations.Range("B12").Value = "TSSR"

$wsAnimations.Columns.Item(1).AutoFit()
$wsAnimations.Columns.Item(2).AutoFit()

# ---------------------------------------------------------------------------
# 5. Activate the last sheet ("Animations"), mirroring the bumped activeTab
# ---------------------------------------------------------------------------
$wsAnimations.Activate()
